$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) values (B2:E2)
$ws.Range("B2").Value = 265.6287409710601
$ws.Range("C2").Value = 290.73128431545058
$ws.Range("D2").Value = 262.21989162037539
$ws.Range("E2").Value = 297.51363840877923

# Row 3 (STR) values (B3:E3)
$ws.Range("B3").Value = 261.16854873030132
$ws.Range("C3").Value = 295.99799930273218
$ws.Range("D3").Value = 261.79447028965973
$ws.Range("E3").Value = 303.44048945715053

# Update the selection to match the new selected range B1:E3
$null = $ws.Range("B1:E3").Select()
